# Actualización automática 2025-07-07 17:30:08
#
# Updates sales figures for client "APOLO CHAMBA KATHERINE YUELISE"
# (row 7) across the three report sheets, plus the dependent
# subtotal/total rows that are derived from that client's figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
#   Row 7  -> client APOLO CHAMBA KATHERINE YUELISE group sales
#   Row 55 -> "x de 53" completion counters per group
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("H7").Value = 536.4    # INODOROS
$wsGrupo.Range("I7").Value = 181.8    # LAVABOS
$wsGrupo.Range("M7").Value = 1284.08  # PORCELANATO

$wsGrupo.Range("H55").Value = "2 de 53"   # INODOROS counter
$wsGrupo.Range("I55").Value = "1 de 53"   # LAVABOS counter
$wsGrupo.Range("M55").Value = "10 de 53"  # PORCELANATO counter

# ---------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
#   Row 7  -> client APOLO CHAMBA KATHERINE YUELISE, julio sales
#   Row 55 -> julio total
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F7").Value = 2002.28
$wsMensual.Range("F55").Value = 25842.17

# ---------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
#   Row 7  -> grupo INODOROS
#   Row 8  -> grupo LAVABOS
#   Row 16 -> grupo PORCELANATO
#   Row 19 -> TOTAL
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D7").Value = 963
$wsCumpl.Range("E7").Value = 1437
$wsCumpl.Range("F7").Value = 0.40125

$wsCumpl.Range("D8").Value = 181.8
$wsCumpl.Range("E8").Value = 818.2
$wsCumpl.Range("F8").Value = 0.1818

$wsCumpl.Range("D16").Value = 22289.84
$wsCumpl.Range("E16").Value = 29536.62
$wsCumpl.Range("F16").Value = 0.4300860988769057

$wsCumpl.Range("D19").Value = 25842.17
$wsCumpl.Range("E19").Value = 87864.28064517914
$wsCumpl.Range("F19").Value = 0.2272709230951238
